$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E12").Formula = "=-(18*0.0007+0.007*4)"
$ws.Range("E13").Formula = "=-(25*0.0007+0.007*5)"
$ws.Range("E14").Formula = "=-(22*0.0007+0.007*4)"
$ws.Range("E15").Formula = "=-(18*0.0007+0.007*3)"
$ws.Range("E16").Formula = "=-(21*0.0007+0.007*4)"

$ws.Range("L21").Select()
